$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.906.80"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "2.216.92"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "291.80"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.60"
$ws.Range("E6").Value = "  +1.78%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.468"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.52"
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0780"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "50.03"
$ws.Range("E12").Value = "  +5.28%  "
$ws.Range("E13").Value = "  +2.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.44"
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("D15").Value = "2.561.30"
$ws.Range("E15").Value = "  +0.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.76"
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").Value = "2.175.26"
$ws.Range("E17").Value = "  -1.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.731"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "39.867.23"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").Value = "0.0₃0887"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.11"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.74"
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.60"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "237.02"
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.45"
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.14"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.24"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.04"
$ws.Range("E30").Value = "  -6.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.09"
$ws.Range("E31").Value = "  +3.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.88"
$ws.Range("E32").Value = "  -2.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.96"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0711"
$ws.Range("E35").Value = "  -0.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.95"
$ws.Range("E36").Value = "  +5.91%  "
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.72"
$ws.Range("E40").Value = "  +2.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.28"
$ws.Range("E41").Value = "  -3.43%  "
$ws.Range("D42").Value = "2.109.46"
$ws.Range("E42").Value = "  +2.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.72"
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0269"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.92"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.77"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.04"
$ws.Range("E47").Value = "  -2.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.68"
$ws.Range("E48").Value = "  +3.37%  "
$ws.Range("D49").Value = "2.435.18"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.47"
$ws.Range("E50").Value = "  +2.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "88.54"
$ws.Range("E51").Value = "  -0.43%  "
